# Auto-generated Excel COM-interop edit script
# Applies updated market-price derived figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across multiple Leve sheets, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1122.1111
$ws.Range("I28").Value = 944.6875
$ws.Range("J28").Value = 2541.5
$ws.Range("K28").Value = 944.6875
$ws.Range("L28").Value = 2541.5
$ws.Range("M28").Value = -459.6875
$ws.Range("N28").Value = -3511.5

$ws.Range("H33").Value = 908938.3
$ws.Range("I33").Value = 1150441.9
$ws.Range("K33").Value = 1150441.9
$ws.Range("M33").Value = -1150212.9

$ws.Range("H74").Value = 4476.4614
$ws.Range("I74").Value = 3359.6
$ws.Range("K74").Value = 3359.6
$ws.Range("M74").Value = -2423.6

$ws.Range("H77").Value = 4476.4614
$ws.Range("I77").Value = 3359.6
$ws.Range("K77").Value = 16798
$ws.Range("M77").Value = -12118

$ws.Range("H88").Value = 342999.8
$ws.Range("J88").Value = 253749.75
$ws.Range("L88").Value = 253749.75
$ws.Range("N88").Value = -254561.75

$ws.Range("H91").Value = 342999.8
$ws.Range("J91").Value = 253749.75
$ws.Range("L91").Value = 253749.75
$ws.Range("N91").Value = -256557.75

$ws.Range("H135").Value = 761.6875
$ws.Range("I135").Value = 654
$ws.Range("J135").Value = 941.1667
$ws.Range("K135").Value = 5886
$ws.Range("L135").Value = 8470.5003
$ws.Range("M135").Value = -3351
$ws.Range("N135").Value = -13540.5003

$ws.Range("H138").Value = 2414.5874
$ws.Range("J138").Value = 3595.1667
$ws.Range("L138").Value = 10785.5001
$ws.Range("N138").Value = -21065.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 176.45454
$ws.Range("J5").Value = 141.25
$ws.Range("L5").Value = 141.25
$ws.Range("N5").Value = -365.25

$ws.Range("H132").Value = 3631.652
$ws.Range("I132").Value = 3695.2131
$ws.Range("K132").Value = 11085.6393
$ws.Range("M132").Value = -8555.639299999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 176.45454
$ws.Range("J4").Value = 141.25
$ws.Range("L4").Value = 141.25
$ws.Range("N4").Value = -371.25

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H107").Value = 5978.613
$ws.Range("I107").Value = 5475.0527
$ws.Range("K107").Value = 5475.0527
$ws.Range("M107").Value = -3555.0527

$ws.Range("H134").Value = 3769.2222
$ws.Range("I134").Value = 3212.0605
$ws.Range("K134").Value = 9636.181500000001
$ws.Range("M134").Value = -7101.181500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 15000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H134").Value = 2679.5454
$ws.Range("I134").Value = 2331
$ws.Range("K134").Value = 6993
$ws.Range("M134").Value = -4458

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 143671
$ws.Range("J34").Value = 1249.5
$ws.Range("L34").Value = 3748.5
$ws.Range("N34").Value = -3916.5

$ws.Range("H39").Value = 45945.695
$ws.Range("I39").Value = 59679.59
$ws.Range("J39").Value = 7033
$ws.Range("K39").Value = 179038.77
$ws.Range("L39").Value = 21099
$ws.Range("M39").Value = -178744.77
$ws.Range("N39").Value = -21687

$ws.Range("H55").Value = 40380.703
$ws.Range("J55").Value = 4993.2
$ws.Range("L55").Value = 14979.6
$ws.Range("N55").Value = -15333.6

$ws.Range("H68").Value = 1001732.06
$ws.Range("J68").Value = 2175756
$ws.Range("L68").Value = 6527268
$ws.Range("N68").Value = -6528890

$ws.Range("H71").Value = 1001732.06
$ws.Range("J71").Value = 2175756
$ws.Range("L71").Value = 19581804
$ws.Range("N71").Value = -19589916

$ws.Range("H81").Value = 160706.58
$ws.Range("J81").Value = 20824.5
$ws.Range("L81").Value = 62473.5
$ws.Range("N81").Value = -64719.5

$ws.Range("H84").Value = 160706.58
$ws.Range("J84").Value = 20824.5
$ws.Range("L84").Value = 187420.5
$ws.Range("N84").Value = -198652.5

$ws.Range("H99").Value = 21143.428
$ws.Range("I99").Value = 7024
$ws.Range("J99").Value = 23496.666
$ws.Range("K99").Value = 21072
$ws.Range("L99").Value = 70489.99800000001
$ws.Range("M99").Value = -18826
$ws.Range("N99").Value = -74981.99800000001

$ws.Range("H109").Value = 69802.53
$ws.Range("I109").Value = 100706.5
$ws.Range("J109").Value = 7994.6
$ws.Range("K109").Value = 302119.5
$ws.Range("L109").Value = 23983.8
$ws.Range("M109").Value = -301079.5
$ws.Range("N109").Value = -26063.8

$ws.Range("H113").Value = 442.94446
$ws.Range("I113").Value = 313.57144
$ws.Range("J113").Value = 525.2727
$ws.Range("K113").Value = 940.71432
$ws.Range("L113").Value = 1575.8181
$ws.Range("M113").Value = 1229.28568
$ws.Range("N113").Value = -5915.8181

$ws.Range("H120").Value = 12510
$ws.Range("I120").Value = 5020
$ws.Range("K120").Value = 15060
$ws.Range("M120").Value = -10222

$ws.Range("H131").Value = 86992.875
$ws.Range("I131").Value = 46719.25
$ws.Range("K131").Value = 140157.75
$ws.Range("M131").Value = -135117.75

$ws.Range("H137").Value = 4264.0713
$ws.Range("I137").Value = 1807.2
$ws.Range("J137").Value = 5629
$ws.Range("K137").Value = 5421.6
$ws.Range("L137").Value = 16887
$ws.Range("M137").Value = -321.6000000000004
$ws.Range("N137").Value = -27087

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3747
$ws.Range("I122").Value = 3333.8823
$ws.Range("J122").Value = 6088
$ws.Range("K122").Value = 10001.6469
$ws.Range("L122").Value = 18264
$ws.Range("M122").Value = -7551.6469
$ws.Range("N122").Value = -23164

$ws.Range("H140").Value = 99499
$ws.Range("J140").Value = 99499
$ws.Range("L140").Value = 99499
$ws.Range("N140").Value = -109859

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 224.8
$ws.Range("I9").Value = 55
$ws.Range("J9").Value = 338
$ws.Range("K9").Value = 55
$ws.Range("L9").Value = 338
$ws.Range("M9").Value = 169
$ws.Range("N9").Value = -786

$ws.Range("H20").Value = 2044500
$ws.Range("J20").Value = 2505625
$ws.Range("L20").Value = 2505625
$ws.Range("N20").Value = -2506077

$ws.Range("H61").Value = 15169585
$ws.Range("I61").Value = 18520926
$ws.Range("K61").Value = 18520926
$ws.Range("M61").Value = -18520724

$ws.Range("H68").Value = 28655.5
$ws.Range("I68").Value = 3619
$ws.Range("J68").Value = 41173.75
$ws.Range("K68").Value = 3619
$ws.Range("L68").Value = 41173.75
$ws.Range("M68").Value = -2870
$ws.Range("N68").Value = -42671.75

$ws.Range("H71").Value = 28655.5
$ws.Range("I71").Value = 3619
$ws.Range("J71").Value = 41173.75
$ws.Range("K71").Value = 18095
$ws.Range("L71").Value = 205868.75
$ws.Range("M71").Value = -14351
$ws.Range("N71").Value = -213356.75

$ws.Range("H113").Value = 15169585
$ws.Range("I113").Value = 18520926
$ws.Range("K113").Value = 18520926
$ws.Range("M113").Value = -18518756

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 4524.9546
$ws.Range("I136").Value = 2956.6
$ws.Range("K136").Value = 8869.799999999999
$ws.Range("M136").Value = -6319.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 35307.5
$ws.Range("J105").Value = 35307.5
$ws.Range("L105").Value = 35307.5
$ws.Range("N105").Value = -42295.5

$ws.Range("H107").Value = 45515960
$ws.Range("I107").Value = 25861.25
$ws.Range("J107").Value = 71510300
$ws.Range("K107").Value = 77583.75
$ws.Range("L107").Value = 214530900
$ws.Range("M107").Value = -75663.75
$ws.Range("N107").Value = -214534740

$ws.Range("H132").Value = 16627.908
$ws.Range("I132").Value = 16627.908
$ws.Range("K132").Value = 49883.724
$ws.Range("M132").Value = -47353.724

$ws.Range("H136").Value = 6203.769
$ws.Range("I136").Value = 6203.769
$ws.Range("K136").Value = 18611.307
$ws.Range("M136").Value = -16061.307

Write-Output "Applied scheduled Sheets market-data refresh."